$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string values introduced by this edit.
$otherFoundLocations = "Other found locations"
$authors2New = "[Kenneth A.%Egol%NULL%2,    Sanjit R.%Konda%NULL%2,    Mackenzie L.%Bird%NULL%2,    Nicket%Dedhia%NULL%2,    Emma K.%Landes%NULL%2,    Rachel A.%Ranson%NULL%2,    Sara J.%Solasz%NULL%2,    Vinay K.%Aggarwal%NULL%2,    Joseph A.%Bosco%NULL%2,    David L.%Furgiuele%NULL%2,    Abhishek%Ganta%NULL%2,    Jason%Gould%NULL%2,    Thomas R.%Lyon%NULL%2,    Toni M.%McLaurin%NULL%2,    Nirmal C.%Tejwani%NULL%2,    Joseph D.%Zuckerman%NULL%2,    Philipp%Leucht%NULL%2]"
$idPmc = "_PMC"
$authors3New = "[Drake G.%LeBrun%NULL%2,    Maxwell A.%Konnaris%NULL%2,    Gregory C.%Ghahramani%NULL%2,    Ajay%Premkumar%NULL%2,    Chris J.%DeFrancesco%NULL%2,    Jordan A.%Gruskay%NULL%2,    Aleksey%Dvorzhinskiy%NULL%2,    Milan S.%Sandhu%NULL%2,    Elan M.%Goldwyn%NULL%2,    Christopher L.%Mendias%NULL%2,    William M.%Ricci%NULL%2]"
$authors4New = "[Josep Maria%Muñoz Vives%NULL%1,    Montsant%Jornet-Gibert%NULL%2,    Montsant%Jornet-Gibert%NULL%0,    J.%Cámara-Cabrera%NULL%2,    J.%Cámara-Cabrera%NULL%0,    Pedro L.%Esteban%NULL%2,    Pedro L.%Esteban%NULL%0,    Laia%Brunet%NULL%2,    Laia%Brunet%NULL%0,    Luis%Delgado-Flores%NULL%2,    Luis%Delgado-Flores%NULL%0,    P.%Camacho-Carrasco%NULL%2,    P.%Camacho-Carrasco%NULL%0,    P.%Torner%NULL%2,    P.%Torner%NULL%0,    Francesc%Marcano-Fernández%NULL%2,    Francesc%Marcano-Fernández%NULL%0,    NULL%NULL%NULL%0,    NULL%NULL%NULL%0]"
$authors5New = "[Zoe B.%Cheung%zoe.cheung@mountsinai.org%1,    David A.%Forsh%NULL%1]"
$idPmcElsevier = "_PMC_elsevier"

# New column I: "Other found locations"
$ws.Range("I1").Value = $otherFoundLocations
$ws.Range("I2").Value = $idPmc
$ws.Range("I3").Value = $idPmc
$ws.Range("I4").Value = $idPmc
$ws.Range("I5").Value = $idPmcElsevier
# Rows 6/7 are blank ("" ) in the source data. A plain Value="" assignment is a
# no-op in Excel (it never materialises a stored empty-string cell), so force
# the cell into existence with a formula that evaluates to the empty string -
# this is the only way to make Excel persist an explicit "" value here.
$ws.Range("I6").Formula = "=""" + """"
$ws.Range("I7").Formula = "=""" + """"

# Updated Authors column (E) values - re-scraped with extra separator whitespace.
$ws.Range("E2").Value = $authors2New
$ws.Range("E3").Value = $authors3New
$ws.Range("E4").Value = $authors4New
$ws.Range("E5").Value = $authors5New
